# Final parent outcome measurements added
# Populate column D (Post Treatment) with the final measurement values
# for each scale/index row, and leave the selection where the user
# finished entering data (D15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 54
    3  = 60
    4  = 57
    5  = 58
    6  = 56
    7  = 57
    8  = 49
    9  = 53
    10 = 53
    11 = 49
    12 = 50
    13 = 51
    14 = 54
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

$ws.Range("D15").Select()
